$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (Price/Volume columns) are stored as text,
# matching the original inlineStr representation, not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2,4).Value = "305.66"
$ws.Cells.Item(2,5).Value = "-4.79%"

$ws.Cells.Item(3,4).Value = "39.50"
$ws.Cells.Item(3,5).Value = "-8.00%"

$ws.Cells.Item(4,4).Value = "5.030"
$ws.Cells.Item(4,5).Value = "-2.76%"

$ws.Cells.Item(5,4).Value = "0.07668"
$ws.Cells.Item(5,5).Value = "-6.11%"

$ws.Cells.Item(6,4).Value = "4.249"
$ws.Cells.Item(6,5).Value = "-1.76%"

$ws.Cells.Item(7,4).Value = "1.600"
$ws.Cells.Item(7,5).Value = "-10.83%"

$ws.Cells.Item(8,4).Value = "0.8848"
$ws.Cells.Item(8,5).Value = "-6.92%"

$ws.Cells.Item(9,4).Value = "0.09751"
$ws.Cells.Item(9,5).Value = "-12.59%"

$ws.Cells.Item(10,5).Value = "-7.55%"

$ws.Cells.Item(11,4).Value = "0.04456"
$ws.Cells.Item(11,5).Value = "-3.71%"

$ws.Cells.Item(12,4).Value = "0.08915"
$ws.Cells.Item(12,5).Value = "-4.65%"

$ws.Cells.Item(13,4).Value = "0.1054"
$ws.Cells.Item(13,5).Value = "-0.59%"

$ws.Cells.Item(14,4).Value = "0.001266"
$ws.Cells.Item(14,5).Value = "-2.02%"

$ws.Cells.Item(15,2).Value = "CoinExToken"
$ws.Cells.Item(15,3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Cells.Item(15,4).Value = "0.04203"
$ws.Cells.Item(15,5).Value = "0.47%"

$ws.Cells.Item(16,2).Value = "TigerCash"
$ws.Cells.Item(16,3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(16,4).Value = "0.005957"
$ws.Cells.Item(16,5).Value = "0.04%"

$ws.Cells.Item(17,2).Value = "LEO"
$ws.Cells.Item(17,3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(17,4).Value = "3.353"
$ws.Cells.Item(17,5).Value = "-0.27%"

$ws.Cells.Item(18,2).Value = "BTSEToken"
$ws.Cells.Item(18,3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(18,4).Value = "2.462"
$ws.Cells.Item(18,5).Value = "-2.69%"

$ws.Cells.Item(19,2).Value = "BitpandaEcosystemToken"
$ws.Cells.Item(19,3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Cells.Item(19,4).Value = "0.3361"
$ws.Cells.Item(19,5).Value = "-0.09%"

$ws.Cells.Item(20,2).Value = "MCDex"
$ws.Cells.Item(20,3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Cells.Item(20,4).Value = "7.041"
$ws.Cells.Item(20,5).Value = "-5.17%"

$ws.Cells.Item(21,2).Value = "ProBitToken"
$ws.Cells.Item(21,3).Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Cells.Item(21,4).Value = "0.1351"
$ws.Cells.Item(21,5).Value = "-2.79%"

$ws.Cells.Item(22,2).Value = "ZBToken"
$ws.Cells.Item(22,3).Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Cells.Item(22,4).Value = "0.3204"
$ws.Cells.Item(22,5).Value = "22.03%"

$ws.Cells.Item(23,4).Value = "0.001195"
$ws.Cells.Item(23,5).Value = "-4.30%"

$ws.Cells.Item(24,4).Value = "0.004055"
$ws.Cells.Item(24,5).Value = "-6.08%"

$ws.Cells.Item(25,4).Value = "0.0001221"
$ws.Cells.Item(25,5).Value = "9.90%"

$ws.Cells.Item(38,4).Value = "0.02313"
$ws.Cells.Item(38,5).Value = "-11.64%"

$ws.Cells.Item(39,4).Value = "0.05128"
$ws.Cells.Item(39,5).Value = "-7.14%"

$ws.Cells.Item(40,4).Value = "0.007928"
$ws.Cells.Item(40,5).Value = "-0.56%"

$ws.Cells.Item(41,4).Value = "0.1323"
$ws.Cells.Item(41,5).Value = "-4.98%"

$ws.Cells.Item(42,4).Value = "0.006490"
$ws.Cells.Item(42,5).Value = "-1.49%"

$ws.Cells.Item(43,4).Value = "0.001979"
$ws.Cells.Item(43,5).Value = "-6.45%"

$ws.Cells.Item(44,4).Value = "0.008640"
$ws.Cells.Item(44,5).Value = "2.38%"

$ws.Cells.Item(45,4).Value = "0.3033"
$ws.Cells.Item(45,5).Value = "-12.42%"

$ws.Cells.Item(46,4).Value = "0.00006535"
$ws.Cells.Item(46,5).Value = "-6.46%"

$ws.Cells.Item(47,4).Value = "0.00000000751"
$ws.Cells.Item(47,5).Value = "0.01%"

$ws.Cells.Item(48,4).Value = "0.007010"
$ws.Cells.Item(48,5).Value = "98.44%"

$ws.Cells.Item(49,4).Value = "0.003391"
$ws.Cells.Item(49,5).Value = "-2.59%"

$ws.Cells.Item(50,4).Value = "0.00002102"
$ws.Cells.Item(50,5).Value = "0.01%"

$ws.Cells.Item(51,4).Value = "0.0002002"
$ws.Cells.Item(51,5).Value = "0.01%"

# Restore default (Normal) style on the price/volume columns so no stray
# number-format styling is left behind on unaffected cells.
$ws.Range("D2:E51").Style = "Normal"
